$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.419.04'
$ws.Range('E2').Value = '  -1.17%  '
$ws.Range('D3').Value = '2.377.28'
$ws.Range('E3').Value = '  +5.71%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '235.79'
$ws.Range('E5').Value = '  +1.47%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.656'
$ws.Range('E6').Value = '  +3.67%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '72.08'
$ws.Range('E7').Value = '  +14.05%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.466'
$ws.Range('E9').Value = '  +4.39%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0973'
$ws.Range('E10').Value = '  -0.83%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '57.22'
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '26.89'
$ws.Range('E12').Value = '  +1.61%  '
$ws.Range('D13').Value = '2.734.55'
$ws.Range('E13').Value = '  +5.84%  '
$ws.Range('E14').Value = '  +0.64%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.91'
$ws.Range('E15').Value = '  +2.57%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.25'
$ws.Range('E16').Value = '  +2.78%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.856'
$ws.Range('E17').Value = '  +3.03%  '
$ws.Range('D18').Value = '2.376.77'
$ws.Range('E18').Value = '  +5.84%  '
$ws.Range('D19').Value = '43.458.04'
$ws.Range('E19').Value = '  -0.87%  '
$ws.Range('E20').Value = '  +0.71%  '
$ws.Range('E21').Value = '  +5.25%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '74.55'
$ws.Range('E22').Value = '  +2.38%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '251.54'
$ws.Range('E23').Value = '  +1.23%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.94'
$ws.Range('E24').Value = '  +18.30%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  +2.16%  '
$ws.Range('E27').Value = '  -0.65%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '23.09'
$ws.Range('E28').Value = '  +9.94%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.01'
$ws.Range('E29').Value = '  +1.85%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '174.65'
$ws.Range('E30').Value = '  +0.74%  '
$ws.Range('E31').Value = '  +8.71%  '
$ws.Range('E32').Value = '  -8.57%  '
$ws.Range('E33').Value = '  +2.08%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.00'
$ws.Range('E34').Value = '  +4.13%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0692'
$ws.Range('E35').Value = '  +1.32%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.07'
$ws.Range('E36').Value = '  +2.75%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.46'
$ws.Range('E37').Value = '  +8.26%  '
$ws.Range('B38').Value = 'THORChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.62'
$ws.Range('E38').Value = '  +3.53%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.67'
$ws.Range('E39').Value = '  +0.28%  '
$ws.Range('E40').Value = '  +0.78%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '18.76'
$ws.Range('E41').Value = '  +10.13%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.93'
$ws.Range('E42').Value = '  +3.61%  '
$ws.Range('B43').Value = 'BinanceUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.18'
$ws.Range('E44').Value = '  +10.17%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '99.85'
$ws.Range('E45').Value = '  +2.02%  '
$ws.Range('E46').Value = '  +2.53%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.48'
$ws.Range('E47').Value = '  +2.64%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0954'
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('D49').Value = '1.453.54'
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('D50').Value = '2.605.59'
$ws.Range('E50').Value = '  +6.02%  '
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.74'
$ws.Range('E51').Value = '  -0.50%  '
